$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: column A's caption changes from "Pressão (MPa)" to "Pressao-Mpa".
# Columns B, C and D keep their original header text untouched.
$ws.Range("A1").Value = "Pressao-Mpa"

# --- Header formatting: the whole header row (A1:D1) is restyled bold + red.
# Fill / border / alignment / number format stay exactly as they were.
$headerRng = $ws.Range("A1:D1")
$headerRng.Font.Bold = $true
$headerRng.Font.Color = 255

# --- Row 1 shrinks now that the header text is short again.
$ws.Rows.Item(1).RowHeight = 15.75

# --- Columns get explicit custom widths (previously a uniform best-fit width).
$ws.Columns.Item(1).ColumnWidth = 31.6658
$ws.Columns.Item(2).ColumnWidth = 34.6603
$ws.Columns.Item(3).ColumnWidth = 48.3306
$ws.Columns.Item(4).ColumnWidth = 49.0055

# --- Selection cursor ends up parked on the newly formatted header range.
$ws.Range("A1:D1").Select() | Out-Null

Write-Host "done"
